$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted right before the existing row 243,
# pushing the former rows 243-342 down to 244-343 (dimension grows by one row).
$ws.Rows.Item(243).Insert()

$ws.Cells.Item(243, 1).Value = 4
$ws.Cells.Item(243, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(243, 3).Value = "Los Lagos"
$ws.Cells.Item(243, 4).Value = 45141
$ws.Cells.Item(243, 5).Value = 10
$ws.Cells.Item(243, 6).Value = 100112009
$ws.Cells.Item(243, 7).Value = "Acelga"
$ws.Cells.Item(243, 8).Value = "Sin especificar"
$ws.Cells.Item(243, 9).Value = "Primera"
$ws.Cells.Item(243, 10).Value = 50
$ws.Cells.Item(243, 11).Value = 10000
$ws.Cells.Item(243, 12).Value = 10000
$ws.Cells.Item(243, 13).Value = 10000
$ws.Cells.Item(243, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(243, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(243, 16).Value = 833
$ws.Cells.Item(243, 17).Value = 12
$ws.Cells.Item(243, 18).Value = "Hortaliza"
